# Weekly crime data update for cs-en-us-044pct.xlsx ("New crime data collected")
# - bump the "Volume 32   Number  NN" report number (24 -> 25)
# - bump the reporting week dates (6/9/2025 - 6/15/2025 -> 6/16/2025 - 6/22/2025)
# - refresh the weekly crime-complaint statistics grid (rows 14-30) with the newly
#   collected week's Week-to-Date / 28-Day / YTD / 2-Year figures and their %Chg columns
#   (a handful of cells flip from the blank "0"/"***.*" placeholder text to real
#   numbers now that this precinct/category has data again, so their number format
#   is restored to match the surrounding integer / one-decimal columns)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# A8 rich text reads "Volume 32   Number  24" -> bump the report number to 25.
$volTxt = $ws.Range("A8").Text
$volPos = $volTxt.IndexOf("24")
$ws.Range("A8").Characters($volPos + 1, 2).Text = "25"

# C9 rich text reads "Report Covering the Week  6/9/2025  Through  6/15/2025"
# -> shift the reporting week forward by 7 days on both ends.
$weekTxt = $ws.Range("C9").Text
$fromPos = $weekTxt.IndexOf("6/9/2025")
$ws.Range("C9").Characters($fromPos + 1, 8).Text = "6/16/2025"

$weekTxt2 = $ws.Range("C9").Text
$thruPos = $weekTxt2.IndexOf("6/15/2025")
$ws.Range("C9").Characters($thruPos + 1, 9).Text = "6/22/2025"

# --- Crime complaint grid updates (rows 14-30) ---------------------------
$ws.Range("L14").Value = -45.454545454545
$ws.Range("N14").Value = -81.25
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 17
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = -22.727272727272
$ws.Range("L15").Value = -19.047619047619
$ws.Range("M15").Value = 6.25
$ws.Range("N15").Value = -58.536585365853
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 36
$ws.Range("G16").Value = 40
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 221
$ws.Range("J16").Value = 253
$ws.Range("K16").Value = -12.648221343873
$ws.Range("L16").Value = -15.969581749049
$ws.Range("M16").Value = -1.339285714285
$ws.Range("N16").Value = -76.687763713080
$ws.Range("C17").Value = 22
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = -4.347826086956
$ws.Range("F17").Value = 87
$ws.Range("G17").Value = 93
$ws.Range("H17").Value = -6.451612903225
$ws.Range("I17").Value = 457
$ws.Range("J17").Value = 513
$ws.Range("K17").Value = -10.916179337232
$ws.Range("L17").Value = -11.605415860735
$ws.Range("M17").Value = 77.821011673151
$ws.Range("N17").Value = -26.409017713365
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -14.285714285714
$ws.Range("F18").Value = 17
$ws.Range("H18").Value = -22.727272727272
$ws.Range("I18").Value = 112
$ws.Range("J18").Value = 137
$ws.Range("K18").Value = -18.248175182481
$ws.Range("L18").Value = -54.285714285714
$ws.Range("M18").Value = -6.666666666666
$ws.Range("N18").Value = -88.465499485066
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 7.692307692307
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = 1.612903225806
$ws.Range("I19").Value = 361
$ws.Range("J19").Value = 422
$ws.Range("K19").Value = -14.454976303317
$ws.Range("L19").Value = 2.849002849002
$ws.Range("M19").Value = 108.670520231214
$ws.Range("N19").Value = -0.824175824175
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -42.857142857142
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -25.925925925925
$ws.Range("I20").Value = 166
$ws.Range("J20").Value = 162
$ws.Range("K20").Value = 2.469135802469
$ws.Range("L20").Value = -28.755364806867
$ws.Range("M20").Value = 78.494623655914
$ws.Range("N20").Value = -75.836972343522
$ws.Range("C21").Value = 59
$ws.Range("D21").Value = 60
$ws.Range("E21").Value = -1.666666666666
$ws.Range("F21").Value = 227
$ws.Range("G21").Value = 245
$ws.Range("H21").Value = -7.346938775510
$ws.Range("I21").Value = 1340
$ws.Range("J21").Value = 1516
$ws.Range("K21").Value = -11.609498680738
$ws.Range("L21").Value = -18.342474101157
$ws.Range("M21").Value = 50.392817059483
$ws.Range("N21").Value = -63.427947598253
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 22
$ws.Range("K22").Value = -24.137931034482
$ws.Range("L22").Value = -8.333333333333
$ws.Range("M22").Value = 4.761904761904
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 25
$ws.Range("J23").Value = 37
$ws.Range("K23").Value = -32.432432432432
$ws.Range("L23").Value = -28.571428571428
$ws.Range("M23").Value = -3.846153846153
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = 5.405405405405
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = 11.811023622047
$ws.Range("I24").Value = 929
$ws.Range("J24").Value = 929
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = -7.1
$ws.Range("M24").Value = 46.529968454258
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 142.857142857143
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 28.571428571428
$ws.Range("I25").Value = 342
$ws.Range("J25").Value = 462
$ws.Range("K25").Value = -25.974025974026
$ws.Range("L25").Value = -35.227272727272
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 23
$ws.Range("E26").Value = -8.695652173913
$ws.Range("F26").Value = 120
$ws.Range("G26").Value = 122
$ws.Range("H26").Value = -1.639344262295
$ws.Range("I26").Value = 593
$ws.Range("J26").Value = 653
$ws.Range("K26").Value = -9.188361408882
$ws.Range("L26").Value = 0.679117147707
$ws.Range("M26").Value = 5.516014234875
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 25
$ws.Range("J27").Value = 34
$ws.Range("K27").Value = -26.470588235294
$ws.Range("L27").Value = -26.470588235294
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = 50
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F28").Value = 14
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 72
$ws.Range("J28").Value = 62
$ws.Range("K28").Value = 16.129032258064
$ws.Range("L28").Value = 5.882352941176
$ws.Range("D29").Value = 2
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J29").Value = 35
$ws.Range("K29").Value = -71.428571428571
$ws.Range("L29").Value = -65.517241379310
$ws.Range("N29").Value = -87.804878048780
$ws.Range("D30").Value = 2
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J30").Value = 24
$ws.Range("K30").Value = -70.833333333333
$ws.Range("L30").Value = -70.833333333333
$ws.Range("N30").Value = -90.277777777777
